# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers: AD1/AE1/AF1, copying the existing header formatting (bold,
#     centered, thin border) from AC1 so the new headers match the rest
#     of the header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Season record values for every player row (2-47): the team finished
#     83-79 with no ties, repeated on every row of the roster.
$ws.Range("AD2:AD47").Value = 83
$ws.Range("AE2:AE47").Value = 79
$ws.Range("AF2:AF47").Value = 0
